$wb = $excel.ActiveWorkbook

# Hyperlink-style blue/underline color used elsewhere in this workbook
# (theme color FF6495ED expressed as BGR int for the COM Font.Color property).
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# zh-cn sheet: row 7 (a8a536f6-...) handback just came in, but it is stale.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhG7 = $wsZh.Range("G7").Value2

$zhI7 = $wsZh.Range("I7")
$zhI7.Value = "a8a536f6-a7af-4647-ab4a-411a8620ffab.md"
$wsZh.Hyperlinks.Add($zhI7, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/50688d82a718d82c67298b23512533006c124ba0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md", "", "", "a8a536f6-a7af-4647-ab4a-411a8620ffab.md") | Out-Null
$zhI7.Font.Underline = 2
$zhI7.Font.Color = $hyperlinkColor

$wsZh.Range("J7").Value = $zhG7
$wsZh.Range("K7").Value = "2016-08-24 08:57:46"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50688d82a718d82c67298b23512533006c124ba0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18cfaf5320d222d7ccb7059798fafc7a84d064c0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md."

# ---------------------------------------------------------------------------
# de-de sheet: same handback-status update for row 7.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deG7 = $wsDe.Range("G7").Value2

$deI7 = $wsDe.Range("I7")
$deI7.Value = "a8a536f6-a7af-4647-ab4a-411a8620ffab.md"
$wsDe.Hyperlinks.Add($deI7, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/50688d82a718d82c67298b23512533006c124ba0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md", "", "", "a8a536f6-a7af-4647-ab4a-411a8620ffab.md") | Out-Null
$deI7.Font.Underline = 2
$deI7.Font.Color = $hyperlinkColor

$wsDe.Range("J7").Value = $deG7
$wsDe.Range("K7").Value = "2016-08-24 08:57:53"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50688d82a718d82c67298b23512533006c124ba0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/18cfaf5320d222d7ccb7059798fafc7a84d064c0/e2e/a8a536f6-a7af-4647-ab4a-411a8620ffab.md."
